$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two data rows "RM 232" (row 26) and "SC 92" (row 28) were removed from
# the dataset entirely. Deleting the higher-numbered row first keeps the
# lower row index ("26") valid for the second delete.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# After the two rows above are removed, the remaining rows shift up to fill
# the gaps. On top of that shift, a handful of "missing value" cells moved
# around within the surviving rows (this file is part of a missing-data /
# imputation benchmark, so individual cells are blanked out or filled back
# in independently of the row deletions).

# RM 14 (row 5): F becomes missing
$ws.Range("F5").Value = ""

# RM 38 (row 8): F is filled back in
$ws.Range("F8").Value = 17.05

# RM 81 (row 12): F becomes missing
$ws.Range("F12").Value = ""

# RM 90 (row 14): F is filled back in
$ws.Range("F14").Value = 17.76

# RM 120 (row 18): F becomes missing
$ws.Range("F18").Value = ""

# SC 5 (now row 26 after the shift): D is filled back in
$ws.Range("D26").Value = -13.8

# SC 101 (now row 27 after the shift): D becomes missing
$ws.Range("D27").Value = ""

# SC 232 (now row 33 after the shift): E is filled back in
$ws.Range("E33").Value = -10.7
